$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.132355213165283
$ws.Range("B1").Value = 2.147614240646362
$ws.Range("C1").Value = 2.744754552841187
$ws.Range("D1").Value = 1.558014750480652
$ws.Range("E1").Value = 0.8219135403633118
